$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 3, pushing existing rows 3+ down by one.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3.
$ws.Range("A3").Value = "baseline2"
$ws.Range("C3").Value = 122.81623181398
$ws.Range("B3").Clear()

# Update the selection to match the recorded view state.
$ws.Range("F7").Select()
